$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# 1. Update the hyperlink display text in C2: "PCI Registered Reports" -> "Peer Community In Registered Reports"
$ws.Cells.Item(2, 3).Value = "\href{https://rr.peercommunityin.org/}{Peer Community In Registered Reports}"

# 2. Move the last record block (row 22, columns A:D) down to row 23 (work bottom-up so nothing is overwritten before it's copied)
$ws.Range("A23:D23").Value = $ws.Range("A22:D22").Value()
$ws.Range("A22:D22").ClearContents()

# 3. Shift column E (the journal list, plus the trailing rows) down by one row: old E8:E23 -> new E9:E24
for ($r = 23; $r -ge 8; $r--) {
    $ws.Cells.Item($r + 1, 5).Value = $ws.Cells.Item($r, 5).Value()
}

# 4. Put the new journal entry at the top of the list (E8), matching the updated C2 hyperlink text
$ws.Cells.Item(8, 5).Value = "\href{https://rr.peercommunityin.org/}{Peer Community In Registered Reports}"

# 5. Update the sheet view: selection now on C2, no pinned top-left cell
$ws.Range("C2").Select()
